# Generate Report for Handoff
#
# Stamp a fresh "Latest Handoff Datetime" on every row that was just handed
# off (i.e. every tracked file row except the "In Translation" row and the
# ".localization-config" / "Not to be localized" row), on both the zh-cn
# and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Column D = "Latest Handoff Datetime"; rows 4,6,7,8,9,10 are the files
# that were (re)handed off in this run.
$handoffRows = @(4, 6, 7, 8, 9, 10)

foreach ($r in $handoffRows) {
    $zhcn.Cells.Item($r, 4).Value = "2016-02-22 14:43:00"
}

foreach ($r in $handoffRows) {
    $dede.Cells.Item($r, 4).Value = "2016-02-22 14:43:14"
}
